$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "62.545.66"
Set-TextValue $ws.Range("E2") "  +1.12%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.436.15"
Set-TextValue $ws.Range("E3") "  +1.31%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "566.49"
Set-TextValue $ws.Range("E5") "  +1.15%  "

# Row 6
Set-TextValue $ws.Range("D6") "145.44"
Set-TextValue $ws.Range("E6") "  +2.48%  "

# Row 7
Set-TextValue $ws.Range("E7") "  -0.04%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.532"
Set-TextValue $ws.Range("E8") "  +0.58%  "

# Row 9
Set-TextValue $ws.Range("E9") "  +2.36%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +0.47%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +1.59%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +2.40%  "

# Row 13
Set-TextValue $ws.Range("D13") "26.84"
Set-TextValue $ws.Range("E13") "  +5.84%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.0000180"
Set-TextValue $ws.Range("E14") "  +4.93%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.860.64"
Set-TextValue $ws.Range("E15") "  +0.88%  "

# Row 16
Set-TextValue $ws.Range("D16") "62.354.97"
Set-TextValue $ws.Range("E16") "  +0.67%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.435.59"
Set-TextValue $ws.Range("E17") "  +1.47%  "

# Row 18
Set-TextValue $ws.Range("D18") "11.26"
Set-TextValue $ws.Range("E18") "  +1.01%  "

# Row 19
Set-TextValue $ws.Range("E19") "  +2.92%  "

# Row 20
Set-TextValue $ws.Range("D20") "323.89"
Set-TextValue $ws.Range("E20") "  +1.13%  "

# Row 21
Set-TextValue $ws.Range("D21") "4.17"
Set-TextValue $ws.Range("E21") "  +1.50%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.999"
Set-TextValue $ws.Range("E22") "  -0.03%  "

# Row 23
Set-TextValue $ws.Range("B23") "Litecoin"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D23") "67.29"
Set-TextValue $ws.Range("E23") "  +3.17%  "

# Row 24
Set-TextValue $ws.Range("B24") "SuiNetwork"
Set-TextValue $ws.Range("C24") "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D24") "1.82"
Set-TextValue $ws.Range("E24") "  +5.82%  "

# Row 25
Set-TextValue $ws.Range("D25") "583.86"
Set-TextValue $ws.Range("E25") "  +3.84%  "

# Row 26
Set-TextValue $ws.Range("D26") "8.56"
Set-TextValue $ws.Range("E26") "  -1.22%  "

# Row 27
Set-TextValue $ws.Range("D27") "0.0₃0994"
Set-TextValue $ws.Range("E27") "  +7.24%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.557.62"
Set-TextValue $ws.Range("E28") "  +1.68%  "

# Row 29
Set-TextValue $ws.Range("D29") "8.44"
Set-TextValue $ws.Range("E29") "  +4.03%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -0.35%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.44"
Set-TextValue $ws.Range("E31") "  +4.43%  "

# Row 32
Set-TextValue $ws.Range("E32") "  -0.31%  "

# Row 33
Set-TextValue $ws.Range("E33") "  +1.26%  "

# Row 34
Set-TextValue $ws.Range("E34") "  +0.23%  "

# Row 35
Set-TextValue $ws.Range("E35") "  +2.46%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.999"
Set-TextValue $ws.Range("E36") "  +0.07%  "

# Row 37
Set-TextValue $ws.Range("E37") "  +0.82%  "

# Row 38
Set-TextValue $ws.Range("D38") "18.76"
Set-TextValue $ws.Range("E38") "  +1.67%  "

# Row 39
Set-TextValue $ws.Range("D39") "5.35"
Set-TextValue $ws.Range("E39") "  -1.19%  "

# Row 40
Set-TextValue $ws.Range("B40") "Monero"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D40") "147.67"
Set-TextValue $ws.Range("E40") "  -2.97%  "

# Row 41
Set-TextValue $ws.Range("B41") "Stacks"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "1.82"
Set-TextValue $ws.Range("E41") "  +2.13%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.999"
Set-TextValue $ws.Range("E42") "  +0.04%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.43"
Set-TextValue $ws.Range("E43") "  +9.17%  "

# Row 44
Set-TextValue $ws.Range("D44") "148.72"
Set-TextValue $ws.Range("E44") "  +0.73%  "

# Row 45
Set-TextValue $ws.Range("E45") "  +2.65%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.0534"
Set-TextValue $ws.Range("E46") "  +1.75%  "

# Row 47
Set-TextValue $ws.Range("D47") "20.50"
Set-TextValue $ws.Range("E47") "  +3.65%  "

# Row 48
Set-TextValue $ws.Range("E48") "  +2.39%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0231"
Set-TextValue $ws.Range("E49") "  +2.93%  "

# Row 50
Set-TextValue $ws.Range("E50") "  +0.83%  "

# Row 51
Set-TextValue $ws.Range("E51") "  +4.82%  "
